$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data point (2026/02/02, 月, 13, 177) was recorded. It sorts chronologically
# before the existing "2026/12/29" block, so insert a new row at 768, which pushes
# the old rows 768-809 down to 769-810 (dimension becomes A1:D810).
$ws.Rows.Item(768).Insert()

# Column A stores dates as plain literal text (e.g. "2026/12/29"), not real Excel
# dates -- every other cell in the column is t="inlineStr". A bare
# Value = "2026/02/02" would be auto-parsed by Excel into a date serial number,
# so write it with a leading apostrophe to force literal text, matching the
# rest of the column.
$ws.Cells.Item(768, 1).Value = "'2026/02/02"
$ws.Cells.Item(768, 2).Value = "月"
$ws.Cells.Item(768, 3).Value = 13
$ws.Cells.Item(768, 4).Value = 177

# The apostrophe-prefix trick leaves a "quote prefix" cell style behind on A768
# that the source row didn't have. Copy the (identically-default) format from
# the row below back onto A768 so the new row's formatting matches its
# neighbours exactly, with only the values differing.
$ws.Cells.Item(769, 1).Copy()
$ws.Cells.Item(768, 1).PasteSpecial(-4122)
